$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg"
$ws.Range("G7").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G14").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G15").Value = "dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G17").Value = "youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G18").Value = "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G21").Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg"
$ws.Range("G24").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G31").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G32").Value = "dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G34").Value = "youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"
$ws.Range("G35").Value = "neveen.nashaat@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
$ws.Range("G37").Value = "Kerelos.zareef@med.asu.edu.eg, System"
$ws.Range("G45").Value = "backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G48").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G51").Value = "eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G52").Value = "abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("G54").Value = "Kerelos.zareef@med.asu.edu.eg, System"
$ws.Range("G62").Value = "backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G65").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G68").Value = "eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G69").Value = "abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg"
$ws.Range("G72").Value = "Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G73").Value = "mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg"
$ws.Range("G77").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G83").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G85").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G86").Value = "Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G92").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg"
$ws.Range("G98").Value = "nourhanhosni@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G102").Value = "neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G103").Value = "neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G109").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G117").Value = "enas.omran@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg`n, takyas@med.asu.edu.eg"
$ws.Range("G119").Value = "neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G120").Value = "neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G126").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G134").Value = "enas.omran@med.asu.edu.eg, user@user.com, Amr-Saeed@med.asu.edu.eg`n, takyas@med.asu.edu.eg"
$ws.Range("G136").Value = "neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, marinasorial@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G137").Value = "neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G143").Value = "mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G144").Value = "mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg"
$ws.Range("G150").Value = "Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G151").Value = "System, marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
